# Update the code with pure ppt
$p = $ppt.ActivePresentation

# Remove the second and third slides (delete from the end first so
# indices of the remaining slides are not shifted mid-operation).
$p.Slides.Item(3).Delete()
$p.Slides.Item(2).Delete()

# Rename the title on the first (now only) slide.
$p.Slides.Item(1).Shapes.Item(1).TextFrame.TextRange.Text = "ChatGPT in Medicine"
